# Update gh-pages to output generated at 1f05065
#
# Schema change applied uniformly to every worksheet:
#   - the old column H "是否有舞台（字符串匹配）" (boolean) is dropped
#   - the old column I "Link" shifts left into H
#   - the old column J "Cover" shifts left into I
# Deleting the old H column does all three in one step and keeps each
# sheet's `dimension` (J -> I) in sync automatically.
#
# On top of that, sheets "展览" (#1) and "全部类型" (#4) -- the two sheets
# that carry the 18 scraped rows -- get refreshed "想去人数" (F) and
# "最低票价" (G) values from the new scrape; G also becomes a genuine
# number (the old unsellable/sold-out text markers collapse to 0).

$wb = $excel.ActiveWorkbook

# New F ("想去人数") / G ("最低票价") values per data row, shared by both
# data-bearing sheets.
$newData = @{
    2  = @{ F = 1431; G = 0 }
    3  = @{ F = 7656; G = 65 }
    4  = @{ F = 531;  G = 0 }
    5  = @{ F = 326;  G = 258 }
    6  = @{ F = 32;   G = 55 }
    7  = @{ F = 22;   G = 55 }
    8  = @{ F = 25;   G = 55 }
    9  = @{ F = 5831; G = 65 }
    10 = @{ F = 149;  G = 168 }
    11 = @{ F = 12;   G = 60 }
    12 = @{ F = 24;   G = 70 }
    13 = @{ F = 1782; G = 39.9 }
    14 = @{ F = 1296; G = 65 }
    15 = @{ F = 279;  G = 68 }
    16 = @{ F = 134;  G = 55 }
    17 = @{ F = 16;   G = 0 }
    18 = @{ F = 5520; G = 60 }
    19 = @{ F = 66;   G = 60 }
}

# Sheet indices (1-based, per Worksheets.Item) that hold the 18 scraped
# data rows; the other sheets ("演出", "本地生活") only have the header row.
$dataSheetIndexes = @(1, 4)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Drop the old boolean column; Link/Cover slide left into H/I.
    $ws.Columns("H").Delete()

    if ($dataSheetIndexes -contains $i) {
        foreach ($row in $newData.Keys) {
            $vals = $newData[$row]
            $ws.Range("F$row").Value = $vals.F
            $ws.Range("G$row").Value = $vals.G
        }
    }
}
